$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.813.77"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "1.731.64"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2748"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.21"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06100"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("D11").Value = "1.734.87"
$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6317"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.497"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "25.810.05"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006610"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").Value = "1.954.83"
$ws.Range("E22").Value = "  -1.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.142"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.728"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.130"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.95%  "

$ws.Range("E26").Value = "  +2.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.497"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("E29").Value = "  -2.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "101.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08279"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.666"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.468"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04472"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.612"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9704"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6118"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.655"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01570"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.919"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3797"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.989"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7163"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05351"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1121"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.169"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.85"
$ws.Range("D50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.584"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.47%  "

Write-Output "Update complete"